# Rename the "AnalysedData" worksheet to "AnalysisSourceData"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AnalysedData")
$ws.Name = "AnalysisSourceData"
